# #5: property boat&car done
# Rebuild the "汽車" (car) sheet (sheet index 3, 1-based) so that row 1 holds
# real field-name headers (like the other property sheets) and rows 2-3 gain
# the common trailing metadata columns (property_category .. index).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Row 1: replace the old (accidentally-duplicated-data) header with the
#     real field names, matching the other sheets' B:G header + new H:N ---
$ws.Cells.Item(1,2).Value = "name"
$ws.Cells.Item(1,3).Value = "capacity"
$ws.Cells.Item(1,4).Value = "owner"
$ws.Cells.Item(1,5).Value = "register_date"
$ws.Cells.Item(1,6).Value = "register_reason"
$ws.Cells.Item(1,7).Value = "acquire_value"
$ws.Cells.Item(1,8).Value = "property_category"
$ws.Cells.Item(1,9).Value = "category"
$ws.Cells.Item(1,10).Value = "date"
$ws.Cells.Item(1,11).Value = "legislator_name"
$ws.Cells.Item(1,12).Value = "legislator_id"
$ws.Cells.Item(1,13).Value = "source_file"
$ws.Cells.Item(1,14).Value = "index"

for ($c = 8; $c -le 14; $c++) {
    $ws.Cells.Item(1,$c).Font.Bold = $true
}

# --- Row 2 (car #29, BENZ) -- existing B:G values are unchanged; add the
#     common trailing metadata columns H:N ---
$ws.Cells.Item(2,8).Value = "land"
$ws.Cells.Item(2,9).Value = "normal"
$ws.Cells.Item(2,10).Value = "2012-04-25"
$ws.Cells.Item(2,11).Value = "詹凱臣"
$ws.Cells.Item(2,12).Value = 1760
$ws.Cells.Item(2,13).Value = "tmp4bd81"
$ws.Cells.Item(2,14).Value = 29

# --- Row 3 (car #30, LEXUS) -- existing B:G values are unchanged; add the
#     common trailing metadata columns H:N ---
$ws.Cells.Item(3,8).Value = "land"
$ws.Cells.Item(3,9).Value = "normal"
$ws.Cells.Item(3,10).Value = "2012-04-25"
$ws.Cells.Item(3,11).Value = "詹凱臣"
$ws.Cells.Item(3,12).Value = 1760
$ws.Cells.Item(3,13).Value = "tmp4bd81"
$ws.Cells.Item(3,14).Value = 30

Write-Output "sheet3 (汽車) rebuilt: A1:N3"
